# Apply the edits described by the diff to List_DPT.xlsx
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)  # "rice"
$ws2 = $wb.Worksheets.Item(2)  # "wheat"

# --- Sheet 1 ("rice") ---
# Remove the second data row (row 3: AJL / Punjab / AY / UP / Rice / 1 / 1095.5)
$ws1.Range("A3:A3").EntireRow.Delete()

# Remove the "Cost" column (column H) entirely
$ws1.Columns.Item(8).Delete()

# Update remaining row 2 values
$ws1.Range("B2").Value = "CWHN"
$ws1.Range("D2").Value = "CWHN"
$ws1.Range("E2").Value = "Punjab"

# --- Sheet 2 ("wheat") ---
# Remove all data rows (rows 2, 3, 4), keeping only the header row
$ws2.Range("A2:A4").EntireRow.Delete()

# Remove the "Cost" column (column H) entirely
$ws2.Columns.Item(8).Delete()

# --- Workbook view ---
# Make the "rice" sheet the active tab (activeTab 1 -> 0)
$ws1.Activate()
